# Applies cryptos list price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text first (matching the original inlineStr cell type), otherwise
# Excel auto-converts them (e.g. "0.0620" -> 6.2E-2) and loses the
# original textual formatting (trailing zeros, etc).
$textForceCells = @(
    "D5"
    "D6"
    "D8"
    "D10"
    "D11"
    "D15"
    "D16"
    "D17"
    "D20"
    "D22"
    "D23"
    "D25"
    "D26"
    "D27"
    "D31"
    "D37"
    "D38"
    "D39"
    "D40"
    "D41"
    "D43"
    "D46"
    "D47"
    "D48"
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Plain numeric-looking price updates (kept as text).
$ws.Range("D5").Value = "214.74"
$ws.Range("D6").Value = "0.516"
$ws.Range("D8").Value = "23.58"
$ws.Range("D10").Value = "0.0620"
$ws.Range("D11").Value = "0.0879"
$ws.Range("D15").Value = "0.551"
$ws.Range("D16").Value = "65.86"
$ws.Range("D17").Value = "251.38"
$ws.Range("D20").Value = "7.56"
$ws.Range("D22").Value = "4.48"
$ws.Range("D23").Value = "9.38"
$ws.Range("D25").Value = "146.89"
$ws.Range("D26").Value = "7.23"
$ws.Range("D27").Value = "16.23"
$ws.Range("D31").Value = "0.0499"
$ws.Range("D37").Value = "0.931"
$ws.Range("D38").Value = "0.580"
$ws.Range("D39").Value = "0.0170"
$ws.Range("D40").Value = "1.03"
$ws.Range("D41").Value = "69.50"
$ws.Range("D43").Value = "5.43"
$ws.Range("D46").Value = "0.791"
$ws.Range("D47").Value = "1.72"
$ws.Range("D48").Value = "88.71"
$ws.Range("D51").Value = "7.78"

# Restore default (Normal) cell style so only the value changed,
# matching the unchanged formatting in the source workbook.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining price/volume/name/link text updates.
$ws.Range("D2").Value = "27.882.46"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "1.667.92"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "1.901.56"
$ws.Range("D13").Value = "1.660.41"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("E17").Value = "  +6.95%  "
$ws.Range("D18").Value = "27.842.30"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  +5.76%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").Value = "1.413.40"
$ws.Range("E34").Value = "  -8.25%  "
$ws.Range("E35").Value = "  -5.74%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -6.05%  "
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("D45").Value = "1.810.75"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").Value = "  +5.21%  "
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E51").Value = "  -5.49%  "
